# "Generate Report for Handoff" — refresh the "Latest Handoff Datetime"
# column (D) for every row that was just (re-)handed off, on both the
# zh-cn and de-de localization-status sheets.
#
# Rows 7, 10, 11, 12, 13, 14, 15, 16 are the files included in this
# handoff batch; they all receive the batch's handoff timestamp
# (which differs per target language).

$wb = $excel.ActiveWorkbook

$handoffRows = @(7, 10, 11, 12, 13, 14, 15, 16)

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($row in $handoffRows) {
    $wsZh.Cells.Item($row, 4).Value = "2016-03-10 08:44:36"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($row in $handoffRows) {
    $wsDe.Cells.Item($row, 4).Value = "2016-03-10 08:44:45"
}
